# Weekly update: insert a new record at row 43 (shifting the existing
# rows 43-90 down to 44-91) for "Jengibre" at Vega Central Mapocho de
# Santiago.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 43:90 down one row, opening up a blank row 43.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new weekly record.
$ws.Range("A43").Value() = 9
$ws.Range("B43").Value() = "Vega Central Mapocho de Santiago"
$ws.Range("C43").Value() = "Metropolitana"
$ws.Range("D43").Value() = 44634
$ws.Range("E43").Value() = 13
$ws.Range("F43").Value() = 100114007
$ws.Range("G43").Value() = "Jengibre"
$ws.Range("H43").Value() = "Sin especificar"
$ws.Range("I43").Value() = "Primera"
$ws.Range("J43").Value() = 520
$ws.Range("K43").Value() = 16000
$ws.Range("L43").Value() = 17000
$ws.Range("M43").Value() = 16500
$ws.Range("N43").Value() = "$/caja 13 kilos"
$ws.Range("O43").Value() = "Perú"
$ws.Range("P43").Value() = 1269
$ws.Range("Q43").Value() = 13
$ws.Range("R43").Value() = "Hortaliza"
